$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column G to fit the new "background concept" content ---
# (target stored width is 33.25 "character" units; this engine quantises
# ColumnWidth to 1/7-wide pixel steps, so 32.5714286 is the closest input
# that lands on the nearest representable stored width, ~33.2857.)
$ws.Columns.Item(7).ColumnWidth = 32.5714286

# --- Move the Worst/Best-case mini table (cols I:K) down by one row ---
# Process bottom-up so we never clobber a value before it has been read.

# New row 10 <- old row 9 (Supplies left)
$ws.Range("I10").Value = "Supplies left"
$ws.Range("J10").Formula = "=J7-J8+J9"
$ws.Range("K10").Formula = "=K7-K8+K9"

# New row 9 <- old row 8 (Total Gain)
$ws.Range("I9").Value = "Total Gain"
$ws.Range("J9").Formula = "=F18"
$ws.Range("K9").Value = 22

# New row 8 <- old row 7 (Total Sink)
$ws.Range("I8").Value = "Total Sink"
$ws.Range("J8").Formula = "=SUM(B18,E18)"
$ws.Range("K8").Formula = "=J8-21"

# New row 7 <- old row 6 (Starting Supplies)
$ws.Range("I7").Value = "Starting Supplies"
$ws.Range("J7").Value = 40
$ws.Range("K7").Formula = "=J7"

# New row 6 <- old row 5 headers (Worst Case / Best Case); old I6 no longer used
$ws.Range("I6").ClearContents()
$ws.Range("J6").Value = "Worst Case"
$ws.Range("K6").Value = "Best Case"

# Old row 5 headers are gone now (table shifted down into row 6)
$ws.Range("J5:K5").ClearContents()

$wb.Save()
